# BURN DOWN CHART SPRINT 3.xlsx - Sprint progress update
#
# On "Hoja2" (the daily effort log), day 4 (row 7) gets logged effort:
# JOSE (C7) and CAMILA (D7) each record 1 hour of work. Every downstream
# formula (H/I/J columns on Hoja2, the L column + chart series on Hoja1)
# recalculates automatically from this single data entry.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

$ws2.Range("C7").Value = 1
$ws2.Range("D7").Value = 1

# The workbook was left open on Hoja2 with D20 selected before the final
# save put the focus back on Hoja1 (which stays the selected tab).
$ws2.Activate() | Out-Null
$ws2.Range("D20").Select() | Out-Null
$ws1.Activate() | Out-Null
